$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.510.32'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '3.151.95'
$ws.Range('E3').Value = '  +0.69%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.149.19'
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.93'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.43%  '
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E13').Value = '  -2.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.98%  '
$ws.Range('D15').Value = '3.673.76'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('E16').Value = '  -1.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = '64.256.43'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('D19').Value = '3.152.33'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.34'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.69'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.04%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.80%  '
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('D35').Value = '0.0₃0845'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.25'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.06%  '
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '51.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '458.68'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.27'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.298'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.96%  '
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').Value = '2.943.97'
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.78'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.05%  '
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.53'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('E50').Value = '  +2.50%  '
$ws.Range('E51').Value = '  -0.63%  '
